# Insert a new weekly record for "Vega Monumental Concepción - Pepino ensalada"
# at row 210, pushing the existing rows 210-244 down to 211-245 (new dimension A1:R245).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 210..244 down by inserting a new blank row at 210.
$ws.Rows(210).Insert()

# After the insert, the data that used to live in row 210 is now in row 211 -
# reuse its constant columns (market/region/product/etc.) for the new row,
# and fill in the new week's own figures.
$ws.Range("A210").Value = $ws.Range("A211").Value2
$ws.Range("B210").Value = $ws.Range("B211").Value2
$ws.Range("C210").Value = $ws.Range("C211").Value2
$ws.Range("D210").Value = 45209
$ws.Range("E210").Value = $ws.Range("E211").Value2
$ws.Range("F210").Value = $ws.Range("F211").Value2
$ws.Range("G210").Value = $ws.Range("G211").Value2
$ws.Range("H210").Value = $ws.Range("H211").Value2
$ws.Range("I210").Value = $ws.Range("I211").Value2
$ws.Range("J210").Value = 150
$ws.Range("K210").Value = 15000
$ws.Range("L210").Value = 15000
$ws.Range("M210").Value = 15000
$ws.Range("N210").Value = $ws.Range("N211").Value2
$ws.Range("O210").Value = $ws.Range("O211").Value2
$ws.Range("P210").Value = 250
$ws.Range("Q210").Value = $ws.Range("Q211").Value2
$ws.Range("R210").Value = $ws.Range("R211").Value2
